# Update "parcial" and "acumulado" score cells in the main grades table.
# Table layout (1-indexed via Word COM):
#   col 3 = 1ER. P., col 4 = 2DO. P., col 6 = 1ER. A., col 7 = 2DO. A.
# Rows 2-5 and 7-8 contain the scores that changed; row 6 (module header) has none.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# Row 2: FISICA II            10 10 / 20 20  ->  9 9 / 15 15
Set-CellText $t 2 3 "9"
Set-CellText $t 2 4 "9"
Set-CellText $t 2 6 "15"
Set-CellText $t 2 7 "15"

# Row 3: CIENCIA, TECNOLOGIA, SOCIEDAD Y VALORES   10 10 / 20 20  ->  9 9 / 15 15
Set-CellText $t 3 3 "9"
Set-CellText $t 3 4 "9"
Set-CellText $t 3 6 "15"
Set-CellText $t 3 7 "15"

# Row 4: CALCULO INTEGRAL     10 10 / 25 25  ->  9 9 / 20 20
Set-CellText $t 4 3 "9"
Set-CellText $t 4 4 "9"
Set-CellText $t 4 6 "20"
Set-CellText $t 4 7 "20"

# Row 5: INGLES V             10 10 / 25 25  ->  9 9 / 20 20
Set-CellText $t 5 3 "9"
Set-CellText $t 5 4 "9"
Set-CellText $t 5 6 "20"
Set-CellText $t 5 7 "20"

# Row 7: CONSTRUYE BASES DE DATOS PARA APLICACIONES WEB   10 10  ->  9 9 (acumulado cols unchanged)
Set-CellText $t 7 3 "9"
Set-CellText $t 7 4 "9"

# Row 8: DESARROLLA APLICACIONES WEB CON CONEXION A BASES DE DATOS  10 10  ->  9 9
Set-CellText $t 8 3 "9"
Set-CellText $t 8 4 "9"
